# rnaSample_1476.xlsx -- "fixed harvester column in rnasamples -- holly
# added S.GISH to harvester in bioSamples"
#
# The "harvester" column (B) for every data row (2-7) gets the new value
# "S.GISH" (previously it duplicated the rnaPreparer column's
# "Retrofitted_1476" text). Everything else on the sheet keeps its
# existing content; a couple of cosmetic layout tweaks (row height,
# column B width, and the active selection) came along with the same
# save in the source workbook, so we reproduce those too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content fix: harvester column B, rows 2-7 -> "S.GISH" ---
$ws.Range("B2:B7").Value = "S.GISH"

# --- cosmetic: row 1 / default row height 15 -> 13.8 ---
$ws.Rows.Item(1).RowHeight = 13.8

# --- cosmetic: column B gets its own (slightly wider) column width ---
$ws.Columns.Item(2).ColumnWidth = 7.996666666666667

# --- cosmetic: selection moves from G2:G7 to the whole column B ---
$ws.Range("B:B").Select()
